$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Insert two new slides right after slide 1 ("Title and Content"
#    layout == ppLayoutObject == 16, same layout used by the other
#    content slides in this deck).
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 16)
$s3 = $p.Slides.Add(3, 16)

# --- New slide 2: "What is continuation passing style?" --------------
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "What is continuation passing style?"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "A process to convert recursive procedures into fully tail recursive procedures`rWhy?`rOnce we explicitly represent the continuation as a data structure certain previously impossible tricks become possible`rThis forms the basis of scheme’s continuations which is a very powerful and weird language feature`rThe fact that this conversion is possible illustrates something deep about the interrelationship between iteration and recursion AND eventually the relationship between data and code"
$tr2.Paragraphs(3).IndentLevel = 2
$tr2.Paragraphs(4).IndentLevel = 2
$tr2.Paragraphs(5).IndentLevel = 2

# --- New slide 3: "I expect you to learn this procedure" -------------
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "I expect you to learn this procedure"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "You will have to both understand`rHow to do it`rWhy it works`rI will try to explain, but this is likely something you need to reflect on yourself, and prove to yourself it is correct"
$tr3.Paragraphs(2).IndentLevel = 2
$tr3.Paragraphs(3).IndentLevel = 2

# ---------------------------------------------------------------------
# 2. Refresh the "today" date field cached on the slide master and on
#    every slide layout (datetimeFigureOut placeholder).
# ---------------------------------------------------------------------
$newDate = "10/13/2023"

$masterShapes = $p.SlideMaster.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $sh = $masterShapes.Item($j)
    if ($sh.Name -like "*Date*" -and $sh.HasTextFrame) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Name -like "*Date*" -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

Write-Output "edit complete"
